# Format the transferred_at date cells (column A, rows 2-3) as real Excel
# dates instead of text strings, so the file imports/exports correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells currently hold the text "2020-01-01". Replace with the actual
# Excel date serial number (2020-01-01) and apply a yyyy-mm-dd date format.
$ws.Range("A2").Value = 43831
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("A3").Value = 43831
$ws.Range("A3").NumberFormat = "yyyy\-mm\-dd"

# Move the active selection, matching the recorded cursor position.
$ws.Range("D25").Select()
